# Update the Workflows sheet:
#  - remove rows for "SIMCE Matemáticas", "DIA Matemáticas", "DIA Lenguaje",
#    "Cálculo Veloz", "Fluidez Lectora", "En Pullinque Todos Leemos", "PDL"
#  - keep "SIMCE Lenguaje" (row becomes row 2) and refresh its output/last_run
#  - keep "DIA Extraer respuestas correctas" (row becomes row 3)
#  - append a new row "SIMCE Lenguaje (Copia)" as row 4 with an empty last_run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows we no longer need, from bottom to top so row indices
# of rows still to be deleted don't shift under us.
$ws.Rows.Item(9).Delete()  # PDL
$ws.Rows.Item(8).Delete()  # En Pullinque Todos Leemos
$ws.Rows.Item(7).Delete()  # Fluidez Lectora
$ws.Rows.Item(6).Delete()  # Cálculo Veloz
$ws.Rows.Item(5).Delete()  # DIA Lenguaje
$ws.Rows.Item(4).Delete()  # DIA Matemáticas
$ws.Rows.Item(2).Delete()  # SIMCE Matemáticas

# After the deletes the sheet now has:
#   row1 header
#   row2 SIMCE Lenguaje   (was row3)
#   row3 DIA Extraer respuestas correctas (was row10)

# Refresh the "SIMCE Lenguaje" row output/last_run values
$ws.Range("D2").Value = "XLSX"
$ws.Range("E2").Value = "2026-01-30 10:33:11"

# Append the new "SIMCE Lenguaje (Copia)" row
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "SIMCE Lenguaje (Copia)"
$ws.Range("C4").Value = "Workflow SIMCE"
$ws.Range("D4").Value = "XLSX"
$ws.Range("E4").Value = ""
